$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new topic "Plaintext Manipulation" in row 24, which pushes
# "Data Science: Questions" and "Data Science: Backwards Design" down by
# one row each, and removes "Data Science: Backwards Design II" (the
# old row 26 content) so the rows below it keep their original topics.
$ws.Range("C24").Value = "Plaintext Manipulation"
$ws.Range("C25").Value = "Data Science: Questions"
$ws.Range("C26").Value = "Data Science: Backwards Design"

# Update the active selection to match the author's final cursor position.
$ws.Range("C1").Select()
